# Issue 404 "Coordinator" -> "Instructor" terminology update, stage 3.
#
# 1. Update the cached "datetimeFigureOut" date-placeholder text on the
#    slide master and every slide layout from 7/8/2012 -> 12/1/2012.
# 2. Rename the "CoordData" entity-box labels on the slide to
#    "InstructorData".

$p = $ppt.ActivePresentation

$oldDate = "7/8/2012"
$newDate = "12/1/2012"
$oldLabel = "CoordData"
$newLabel = "InstructorData"

# --- 1. Date placeholder on the slide master ---------------------------
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- 1b. Date placeholder on every slide layout -------------------------
$layouts = $master.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    $layout = $layouts.Item($L)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# --- 2. CoordData -> InstructorData on slide 1 --------------------------
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $sh = $slide.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        if ($sh.TextFrame.TextRange.Text -eq $oldLabel) {
            $sh.TextFrame.TextRange.Text = $newLabel
        }
    }
}
